$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.573935475788744
$ws.Range("C2").Value = 0.2129577057665131
$ws.Range("D2").Value = 0.07831156368222025
$ws.Range("E2").Value = 0.08597221369884878
$ws.Range("G2").Value = 0.002457871387310795
$ws.Range("L2").Value = 0.2274027238441789
$ws.Range("O2").Value = 3.794447473458831
$ws.Range("B3").Value = 1.442913115761371
$ws.Range("C3").Value = 0.1966891919050511
$ws.Range("D3").Value = 0.07105671735978092
$ws.Range("E3").Value = 0.08681787219266823
$ws.Range("G3").Value = 0.002461721217078994
$ws.Range("L3").Value = 0.2175416844065694
$ws.Range("O3").Value = 3.798199984654332
$ws.Range("B4").Value = 1.362781786398614
$ws.Range("C4").Value = 0.1866396841446658
$ws.Range("D4").Value = 0.06664004036537108
$ws.Range("E4").Value = 0.0873850738890809
$ws.Range("G4").Value = 0.002464210237348917
$ws.Range("L4").Value = 0.2116030452407784
$ws.Range("O4").Value = 3.803550902884069
$ws.Range("B5").Value = 1.330208450467865
$ws.Range("C5").Value = 0.182529403982727
$ws.Range("D5").Value = 0.06484964467021825
$ws.Range("E5").Value = 0.08762826724217909
$ws.Range("G5").Value = 0.002465256118770877
$ws.Range("L5").Value = 0.2092121862673366
$ws.Range("O5").Value = 3.806495172706008
$ws.Range("B6").Value = 1.324804589006419
$ws.Range("C6").Value = 0.1818459938768342
$ws.Range("D6").Value = 0.06455291853411893
$ws.Range("E6").Value = 0.08766937714780632
$ws.Range("G6").Value = 0.002465431697266827
$ws.Range("L6").Value = 0.2088169485714104
$ws.Range("O6").Value = 3.807030122591499
$ws.Range("B7").Value = 1.362342162187304
$ws.Range("C7").Value = 0.1865843120214095
$ws.Range("D7").Value = 0.06661585635589518
$ws.Range("E7").Value = 0.08738830488236005
$ws.Range("G7").Value = 0.002464224214444256
$ws.Range("L7").Value = 0.2115706831154966
$ws.Range("O7").Value = 3.803587521279468
$ws.Range("B8").Value = 1.528693818795944
$ws.Range("C8").Value = 0.2073610130329939
$ws.Range("D8").Value = 0.07580219328406201
$ws.Range("E8").Value = 0.08625383770445438
$ws.Range("G8").Value = 0.002459172881915786
$ws.Range("L8").Value = 0.2239785385532826
$ws.Range("O8").Value = 3.795107489445314
$ws.Range("B9").Value = 1.857392059153028
$ws.Range("C9").Value = 0.2476163619250542
$ws.Range("D9").Value = 0.09412111809409396
$ws.Range("E9").Value = 0.08441004193914736
$ws.Range("G9").Value = 0.002450256117376818
$ws.Range("L9").Value = 0.2492331490969093
$ws.Range("O9").Value = 3.802768556042281
$ws.Range("B10").Value = 2.100385250036766
$ws.Range("C10").Value = 0.2768881927411826
$ws.Range("D10").Value = 0.1077731700182909
$ws.Range("E10").Value = 0.08328806037072667
$ws.Range("G10").Value = 0.002444301372142688
$ws.Range("L10").Value = 0.2683551394330408
$ws.Range("O10").Value = 3.823371891331448
$ws.Range("B11").Value = 2.211253501630324
$ws.Range("C11").Value = 0.2901376316416417
$ws.Range("D11").Value = 0.11402747724118
$ws.Range("E11").Value = 0.08282826137288524
$ws.Range("G11").Value = 0.002441720532934594
$ws.Range("L11").Value = 0.2771786524517381
$ws.Range("O11").Value = 3.836032929927228
$ws.Range("B12").Value = 2.253283168298083
$ws.Range("C12").Value = 0.2951451409776951
$ws.Range("D12").Value = 0.1164022433007261
$ws.Range("E12").Value = 0.08266143186523678
$ws.Range("G12").Value = 0.002440761538591396
$ws.Range("L12").Value = 0.2805378893092865
$ws.Range("O12").Value = 3.841303036670581
$ws.Range("B13").Value = 2.244229287162284
$ws.Range("C13").Value = 0.2940671218416639
$ws.Range("D13").Value = 0.1158905092129601
$ws.Range("E13").Value = 0.08269703734885958
$ws.Range("G13").Value = 0.002440967262065339
$ws.Range("L13").Value = 0.2798136183439368
$ws.Range("O13").Value = 3.840146826265709
$ws.Range("B14").Value = 2.214710386588138
$ws.Range("C14").Value = 0.2905497993426138
$ws.Range("D14").Value = 0.1142227220367005
$ws.Range("E14").Value = 0.08281439011703284
$ws.Range("G14").Value = 0.002441641269234245
$ws.Range("L14").Value = 0.2774546587748432
$ws.Range("O14").Value = 3.83645695625313
$ws.Range("B15").Value = 2.196635200545416
$ws.Range("C15").Value = 0.2883940598790957
$ws.Range("D15").Value = 0.113201989518501
$ws.Range("E15").Value = 0.08288722129930903
$ws.Range("G15").Value = 0.002442056500257731
$ws.Range("L15").Value = 0.2760120678277076
$ws.Range("O15").Value = 3.834258829723069
$ws.Range("B16").Value = 2.093146282815724
$ws.Range("C16").Value = 0.2760209554232915
$ws.Range("D16").Value = 0.1073653269458106
$ws.Range("E16").Value = 0.08331912835575217
$ws.Range("G16").Value = 0.002444472603762508
$ws.Range("L16").Value = 0.2677810150721882
$ws.Range("O16").Value = 3.822610871675977
$ws.Range("B17").Value = 2.029742762317085
$ws.Range("C17").Value = 0.2684132840077496
$ws.Range("D17").Value = 0.1037960292880769
$ws.Range("E17").Value = 0.08359705583781896
$ws.Range("G17").Value = 0.002445987523424259
$ws.Range("L17").Value = 0.262763510841225
$ws.Range("O17").Value = 3.81630945192336
$ws.Range("B18").Value = 1.993305805984107
$ws.Range("C18").Value = 0.2640313088974438
$ws.Range("D18").Value = 0.1017471928649059
$ws.Range("E18").Value = 0.08376167406924573
$ws.Range("G18").Value = 0.002446870919023596
$ws.Range("L18").Value = 0.2598893231344306
$ws.Range("O18").Value = 3.812994403242101
$ws.Range("B19").Value = 1.980974249430119
$ws.Range("C19").Value = 0.2625465814703318
$ws.Range("D19").Value = 0.1010541988163709
$ws.Range("E19").Value = 0.08381822847291609
$ws.Range("G19").Value = 0.002447172094694579
$ws.Range("L19").Value = 0.2589181898838291
$ws.Range("O19").Value = 3.811925033629365
$ws.Range("B20").Value = 2.036488967618936
$ws.Range("C20").Value = 0.2692237813441807
$ws.Range("D20").Value = 0.1041755591865154
$ws.Range("E20").Value = 0.08356697708773808
$ws.Range("G20").Value = 0.002445825010666699
$ws.Range("L20").Value = 0.2632964167597578
$ws.Range("O20").Value = 3.816948211030422
$ws.Range("B21").Value = 2.223379556250791
$ws.Range("C21").Value = 0.2915831888937817
$ws.Range("D21").Value = 0.114712417341579
$ws.Range("E21").Value = 0.08277972292796676
$ws.Range("G21").Value = 0.002441442800912055
$ws.Range("L21").Value = 0.2781470547203639
$ws.Range("O21").Value = 3.837527829535759
$ws.Range("B22").Value = 2.345792541646176
$ws.Range("C22").Value = 0.3061393394134484
$ws.Range("D22").Value = 0.1216361794842129
$ws.Range("E22").Value = 0.08230768252279752
$ws.Range("G22").Value = 0.002438685475928252
$ws.Range("L22").Value = 0.2879575209940555
$ws.Range("O22").Value = 3.853751518443119
$ws.Range("B23").Value = 2.280434154751163
$ws.Range("C23").Value = 0.2983757323861482
$ws.Range("D23").Value = 0.1179373963290971
$ws.Range("E23").Value = 0.08255572927987131
$ws.Range("G23").Value = 0.002440147379465309
$ws.Range("L23").Value = 0.28271190370981
$ws.Range("O23").Value = 3.844837919128508
$ws.Range("B24").Value = 2.033438962733442
$ws.Range("C24").Value = 0.2688573811021513
$ws.Range("D24").Value = 0.104003963759979
$ws.Range("E24").Value = 0.08358056063395658
$ws.Range("G24").Value = 0.002445898443866228
$ws.Range("L24").Value = 0.2630554574725465
$ws.Range("O24").Value = 3.816658469676099
$ws.Range("B25").Value = 1.768206326399252
$ws.Range("C25").Value = 0.2367791206585821
$ws.Range("D25").Value = 0.089131988593806
$ws.Range("E25").Value = 0.08486801376890796
$ws.Range("G25").Value = 0.002452563140357036
$ws.Range("L25").Value = 0.2423018181656857
$ws.Range("O25").Value = 3.798078369299276

Write-Host "Updated $($ws.Name) with new pl_mw values (380 kV case)"
